$wb = $excel.ActiveWorkbook

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Ironbark No. 1 Coal Mine, Australia, M0052, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet updates ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

for ($row = 2; $row -le 11; $row++) {
    $data.Cells.Item($row, 19).Value = $newVersion
}
